$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "rewardpos" header and its values in column B (B1:B12),
# leaving only the "rewardval" column (A) populated.
$ws.Range("B1:B12").ClearContents()

# Select B1:B12 to match the resulting active selection in the file.
$ws.Range("B1:B12").Select()
